{"js": "// Replace each two-digit multiplication expression in the document's\n// table cells with its new value, per the commit diff. Every \"old\" text\n// value is unique within the document, so a direct search+replace for\n// each pair is safe and unambiguous.\nconst replacements = [\n  [\"29\u00d724=\", \"63\u00d733=\"],\n  [\"98\u00d794=\", \"48\u00d729=\"],\n  [\"49\u00d756=\", \"43\u00d756=\"],\n  [\"43\u00d755=\", \"50\u00d724=\"],\n  [\"80\u00d725=\", \"98\u00d724=\"],\n  [\"57\u00d758=\", \"82\u00d742=\"],\n  [\"29\u00d750=\", \"79\u00d737=\"],\n  [\"86\u00d732=\", \"70\u00d774=\"],\n  [\"54\u00d721=\", \"45\u00d786=\"],\n  [\"86\u00d745=\", \"96\u00d777=\"],\n  [\"13\u00d717=\", \"85\u00d775=\"],\n  [\"56\u00d728=\", \"23\u00d775=\"],\n  [\"29\u00d778=\", \"83\u00d716=\"],\n  [\"18\u00d784=\", \"40\u00d732=\"],\n  [\"70\u00d716=\", \"14\u00d747=\"],\n  [\"45\u00d799=\", \"70\u00d744=\"],\n  [\"18\u00d745=\", \"31\u00d741=\"],\n  [\"43\u00d778=\", \"53\u00d729=\"],\n  [\"39\u00d797=\", \"89\u00d723=\"],\n  [\"42\u00d756=\", \"77\u00d717=\"],\n  [\"15\u00d713=\", \"33\u00d735=\"],\n  [\"61\u00d746=\", \"19\u00d746=\"],\n  [\"21\u00d796=\", \"37\u00d757=\"],\n  [\"82\u00d719=\", \"26\u00d722=\"],\n  [\"91\u00d753=\", \"17\u00d753=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document's\n# table cells with its new value, per the commit diff. Every \"old\" text\n# value is unique within the document, so a direct Find/Replace for\n# each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"29\u00d724=\", \"63\u00d733=\"),\n    @(\"98\u00d794=\", \"48\u00d729=\"),\n    @(\"49\u00d756=\", \"43\u00d756=\"),\n    @(\"43\u00d755=\", \"50\u00d724=\"),\n    @(\"80\u00d725=\", \"98\u00d724=\"),\n    @(\"57\u00d758=\", \"82\u00d742=\"),\n    @(\"29\u00d750=\", \"79\u00d737=\"),\n    @(\"86\u00d732=\", \"70\u00d774=\"),\n    @(\"54\u00d721=\", \"45\u00d786=\"),\n    @(\"86\u00d745=\", \"96\u00d777=\"),\n    @(\"13\u00d717=\", \"85\u00d775=\"),\n    @(\"56\u00d728=\", \"23\u00d775=\"),\n    @(\"29\u00d778=\", \"83\u00d716=\"),\n    @(\"18\u00d784=\", \"40\u00d732=\"),\n    @(\"70\u00d716=\", \"14\u00d747=\"),\n    @(\"45\u00d799=\", \"70\u00d744=\"),\n    @(\"18\u00d745=\", \"31\u00d741=\"),\n    @(\"43\u00d778=\", \"53\u00d729=\"),\n    @(\"39\u00d797=\", \"89\u00d723=\"),\n    @(\"42\u00d756=\", \"77\u00d717=\"),\n    @(\"15\u00d713=\", \"33\u00d735=\"),\n    @(\"61\u00d746=\", \"19\u00d746=\"),\n    @(\"21\u00d796=\", \"37\u00d757=\"),\n    @(\"82\u00d719=\", \"26\u00d722=\"),\n    @(\"91\u00d753=\", \"17\u00d753=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n"}
